$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.169591333333333
$ws.Cells.Item(2, 8).Value = 3.508774
$ws.Cells.Item(2, 9).Value = 0.05593990076588554
$ws.Cells.Item(2, 10).Value = 0.05593990076588554
$ws.Cells.Item(2, 13).Value = 86.89540866666668
$ws.Cells.Item(2, 14).Value = 260.686226
$ws.Cells.Item(2, 15).Value = 0.319779657009892
$ws.Cells.Item(2, 16).Value = 0.3197796570098919
$ws.Cells.Item(2, 17).Value = 101.6321168829916
$ws.Cells.Item(2, 18).Value = 914.6890519469241
$ws.Cells.Item(2, 19).Value = 0.01788844228008227
$ws.Cells.Item(2, 20).Value = 0.01788844228008227
$ws.Cells.Item(3, 7).Value = 1.169591333333333
$ws.Cells.Item(3, 8).Value = 3.508774
$ws.Cells.Item(3, 9).Value = 0.05593990076588554
$ws.Cells.Item(3, 10).Value = 0.05593990076588554
$ws.Cells.Item(3, 15).Value = 0.1999969065479545
$ws.Cells.Item(3, 16).Value = 0.1999969065479545
$ws.Cells.Item(3, 17).Value = 63.56285816483222
$ws.Cells.Item(3, 18).Value = 572.0657234834899
$ws.Cells.Item(3, 19).Value = 0.01118780710577666
$ws.Cells.Item(3, 20).Value = 0.01118780710577666
$ws.Cells.Item(4, 7).Value = 1.169591333333333
$ws.Cells.Item(4, 8).Value = 3.508774
$ws.Cells.Item(4, 9).Value = 0.05593990076588554
$ws.Cells.Item(4, 10).Value = 0.05593990076588554
$ws.Cells.Item(4, 13).Value = 60.92601633333334
$ws.Cells.Item(4, 14).Value = 182.778049
$ws.Cells.Item(4, 15).Value = 0.224210932487692
$ws.Cells.Item(4, 16).Value = 0.224210932487692
$ws.Cells.Item(4, 17).Value = 71.25854067799177
$ws.Cells.Item(4, 18).Value = 641.326866101926
$ws.Cells.Item(4, 19).Value = 0.01254233731398815
$ws.Cells.Item(4, 20).Value = 0.01254233731398815
$ws.Cells.Item(5, 7).Value = 1.169591333333333
$ws.Cells.Item(5, 8).Value = 3.508774
$ws.Cells.Item(5, 9).Value = 0.05593990076588554
$ws.Cells.Item(5, 10).Value = 0.05593990076588554
$ws.Cells.Item(5, 13).Value = 7.809668333333332
$ws.Cells.Item(5, 14).Value = 23.429005
$ws.Cells.Item(5, 15).Value = 0.02873998867505581
$ws.Cells.Item(5, 16).Value = 0.02873998867505581
$ws.Cells.Item(5, 17).Value = 9.134120398874442
$ws.Cells.Item(5, 18).Value = 82.20708358986998
$ws.Cells.Item(5, 19).Value = 0.001607712114495296
$ws.Cells.Item(5, 20).Value = 0.001607712114495296
$ws.Cells.Item(6, 7).Value = 1.169591333333333
$ws.Cells.Item(6, 8).Value = 3.508774
$ws.Cells.Item(6, 9).Value = 0.05593990076588554
$ws.Cells.Item(6, 10).Value = 0.05593990076588554
$ws.Cells.Item(6, 13).Value = 61.75795633333333
$ws.Cells.Item(6, 14).Value = 185.273869
$ws.Cells.Item(6, 15).Value = 0.2272725152794058
$ws.Cells.Item(6, 16).Value = 0.2272725152794058
$ws.Cells.Item(6, 17).Value = 72.2315704918451
$ws.Cells.Item(6, 18).Value = 650.0841344266059
$ws.Cells.Item(6, 19).Value = 0.01271360195154317
$ws.Cells.Item(6, 20).Value = 0.01271360195154317
$ws.Cells.Item(7, 9).Value = 0.2496787433529823
$ws.Cells.Item(7, 10).Value = 0.2496787433529823
$ws.Cells.Item(7, 13).Value = 86.89540866666668
$ws.Cells.Item(7, 14).Value = 260.686226
$ws.Cells.Item(7, 15).Value = 0.319779657009892
$ws.Cells.Item(7, 16).Value = 0.3197796570098919
$ws.Cells.Item(7, 17).Value = 453.6185956755165
$ws.Cells.Item(7, 18).Value = 4082.567361079648
$ws.Cells.Item(7, 19).Value = 0.07984218291207752
$ws.Cells.Item(7, 20).Value = 0.07984218291207749
$ws.Cells.Item(8, 9).Value = 0.2496787433529823
$ws.Cells.Item(8, 10).Value = 0.2496787433529823
$ws.Cells.Item(8, 15).Value = 0.1999969065479545
$ws.Cells.Item(8, 16).Value = 0.1999969065479545
$ws.Cells.Item(8, 19).Value = 0.04993497630137712
$ws.Cells.Item(8, 20).Value = 0.04993497630137712
$ws.Cells.Item(9, 9).Value = 0.2496787433529823
$ws.Cells.Item(9, 10).Value = 0.2496787433529823
$ws.Cells.Item(9, 13).Value = 60.92601633333334
$ws.Cells.Item(9, 14).Value = 182.778049
$ws.Cells.Item(9, 15).Value = 0.224210932487692
$ws.Cells.Item(9, 16).Value = 0.224210932487692
$ws.Cells.Item(9, 17).Value = 318.0510270139503
$ws.Cells.Item(9, 18).Value = 2862.459243125552
$ws.Cells.Item(9, 19).Value = 0.05598070386952729
$ws.Cells.Item(9, 20).Value = 0.05598070386952728
$ws.Cells.Item(10, 9).Value = 0.2496787433529823
$ws.Cells.Item(10, 10).Value = 0.2496787433529823
$ws.Cells.Item(10, 13).Value = 7.809668333333332
$ws.Cells.Item(10, 14).Value = 23.429005
$ws.Cells.Item(10, 15).Value = 0.02873998867505581
$ws.Cells.Item(10, 16).Value = 0.02873998867505581
$ws.Cells.Item(10, 17).Value = 40.76867623291555
$ws.Cells.Item(10, 18).Value = 366.91808609624
$ws.Cells.Item(10, 19).Value = 0.007175764256366879
$ws.Cells.Item(10, 20).Value = 0.007175764256366878
$ws.Cells.Item(11, 9).Value = 0.2496787433529823
$ws.Cells.Item(11, 10).Value = 0.2496787433529823
$ws.Cells.Item(11, 13).Value = 61.75795633333333
$ws.Cells.Item(11, 14).Value = 185.273869
$ws.Cells.Item(11, 15).Value = 0.2272725152794058
$ws.Cells.Item(11, 16).Value = 0.2272725152794058
$ws.Cells.Item(11, 17).Value = 322.3939889756569
$ws.Cells.Item(11, 18).Value = 2901.545900780912
$ws.Cells.Item(11, 19).Value = 0.05674511601363352
$ws.Cells.Item(11, 20).Value = 0.0567451160136335
$ws.Cells.Item(12, 7).Value = 4.885583666666666
$ws.Cells.Item(12, 8).Value = 14.656751
$ws.Cells.Item(12, 9).Value = 0.2336705631341014
$ws.Cells.Item(12, 10).Value = 0.2336705631341014
$ws.Cells.Item(12, 13).Value = 86.89540866666668
$ws.Cells.Item(12, 14).Value = 260.686226
$ws.Cells.Item(12, 15).Value = 0.319779657009892
$ws.Cells.Item(12, 16).Value = 0.3197796570098919
$ws.Cells.Item(12, 17).Value = 424.5347892901918
$ws.Cells.Item(12, 18).Value = 3820.813103611727
$ws.Cells.Item(12, 19).Value = 0.07472309253233125
$ws.Cells.Item(12, 20).Value = 0.07472309253233125
$ws.Cells.Item(13, 7).Value = 4.885583666666666
$ws.Cells.Item(13, 8).Value = 14.656751
$ws.Cells.Item(13, 9).Value = 0.2336705631341014
$ws.Cells.Item(13, 10).Value = 0.2336705631341014
$ws.Cells.Item(13, 15).Value = 0.1999969065479545
$ws.Cells.Item(13, 16).Value = 0.1999969065479545
$ws.Cells.Item(13, 17).Value = 265.5129640638761
$ws.Cells.Item(13, 18).Value = 2389.616676574885
$ws.Cells.Item(13, 19).Value = 0.0467333897781388
$ws.Cells.Item(13, 20).Value = 0.04673338977813879
$ws.Cells.Item(14, 7).Value = 4.885583666666666
$ws.Cells.Item(14, 8).Value = 14.656751
$ws.Cells.Item(14, 9).Value = 0.2336705631341014
$ws.Cells.Item(14, 10).Value = 0.2336705631341014
$ws.Cells.Item(14, 13).Value = 60.92601633333334
$ws.Cells.Item(14, 14).Value = 182.778049
$ws.Cells.Item(14, 15).Value = 0.224210932487692
$ws.Cells.Item(14, 16).Value = 0.224210932487692
$ws.Cells.Item(14, 17).Value = 297.6591502731999
$ws.Cells.Item(14, 18).Value = 2678.932352458799
$ws.Cells.Item(14, 19).Value = 0.05239149485522099
$ws.Cells.Item(14, 20).Value = 0.05239149485522098
$ws.Cells.Item(15, 7).Value = 4.885583666666666
$ws.Cells.Item(15, 8).Value = 14.656751
$ws.Cells.Item(15, 9).Value = 0.2336705631341014
$ws.Cells.Item(15, 10).Value = 0.2336705631341014
$ws.Cells.Item(15, 13).Value = 7.809668333333332
$ws.Cells.Item(15, 14).Value = 23.429005
$ws.Cells.Item(15, 15).Value = 0.02873998867505581
$ws.Cells.Item(15, 16).Value = 0.02873998867505581
$ws.Cells.Item(15, 17).Value = 38.15478805141721
$ws.Cells.Item(15, 18).Value = 343.393092462755
$ws.Cells.Item(15, 19).Value = 0.00671568933816799
$ws.Cells.Item(15, 20).Value = 0.00671568933816799
$ws.Cells.Item(16, 7).Value = 4.885583666666666
$ws.Cells.Item(16, 8).Value = 14.656751
$ws.Cells.Item(16, 9).Value = 0.2336705631341014
$ws.Cells.Item(16, 10).Value = 0.2336705631341014
$ws.Cells.Item(16, 13).Value = 61.75795633333333
$ws.Cells.Item(16, 14).Value = 185.273869
$ws.Cells.Item(16, 15).Value = 0.2272725152794058
$ws.Cells.Item(16, 16).Value = 0.2272725152794058
$ws.Cells.Item(16, 17).Value = 301.7236627488466
$ws.Cells.Item(16, 18).Value = 2715.512964739619
$ws.Cells.Item(16, 19).Value = 0.05310689663024244
$ws.Cells.Item(16, 20).Value = 0.05310689663024243
$ws.Cells.Item(17, 7).Value = 2.891468666666666
$ws.Cells.Item(17, 8).Value = 8.674405999999999
$ws.Cells.Item(17, 9).Value = 0.1382948604962879
$ws.Cells.Item(17, 10).Value = 0.1382948604962879
$ws.Cells.Item(17, 13).Value = 86.89540866666668
$ws.Cells.Item(17, 14).Value = 260.686226
$ws.Cells.Item(17, 15).Value = 0.319779657009892
$ws.Cells.Item(17, 16).Value = 0.3197796570098919
$ws.Cells.Item(17, 17).Value = 251.2553514368618
$ws.Cells.Item(17, 18).Value = 2261.298162931756
$ws.Cells.Item(17, 19).Value = 0.04422388305573381
$ws.Cells.Item(17, 20).Value = 0.04422388305573379
$ws.Cells.Item(18, 7).Value = 2.891468666666666
$ws.Cells.Item(18, 8).Value = 8.674405999999999
$ws.Cells.Item(18, 9).Value = 0.1382948604962879
$ws.Cells.Item(18, 10).Value = 0.1382948604962879
$ws.Cells.Item(18, 15).Value = 0.1999969065479545
$ws.Cells.Item(18, 16).Value = 0.1999969065479545
$ws.Cells.Item(18, 17).Value = 157.1403681862011
$ws.Cells.Item(18, 18).Value = 1414.26331367581
$ws.Cells.Item(18, 19).Value = 0.02765854429073851
$ws.Cells.Item(18, 20).Value = 0.0276585442907385
$ws.Cells.Item(19, 7).Value = 2.891468666666666
$ws.Cells.Item(19, 8).Value = 8.674405999999999
$ws.Cells.Item(19, 9).Value = 0.1382948604962879
$ws.Cells.Item(19, 10).Value = 0.1382948604962879
$ws.Cells.Item(19, 13).Value = 60.92601633333334
$ws.Cells.Item(19, 14).Value = 182.778049
$ws.Cells.Item(19, 15).Value = 0.224210932487692
$ws.Cells.Item(19, 16).Value = 0.224210932487692
$ws.Cells.Item(19, 17).Value = 176.1656672126549
$ws.Cells.Item(19, 18).Value = 1585.491004913894
$ws.Cells.Item(19, 19).Value = 0.03100721963012799
$ws.Cells.Item(19, 20).Value = 0.03100721963012798
$ws.Cells.Item(20, 7).Value = 2.891468666666666
$ws.Cells.Item(20, 8).Value = 8.674405999999999
$ws.Cells.Item(20, 9).Value = 0.1382948604962879
$ws.Cells.Item(20, 10).Value = 0.1382948604962879
$ws.Cells.Item(20, 13).Value = 7.809668333333332
$ws.Cells.Item(20, 14).Value = 23.429005
$ws.Cells.Item(20, 15).Value = 0.02873998867505581
$ws.Cells.Item(20, 16).Value = 0.02873998867505581
$ws.Cells.Item(20, 17).Value = 22.58141128289222
$ws.Cells.Item(20, 18).Value = 203.2327015460299
$ws.Cells.Item(20, 19).Value = 0.003974592724481738
$ws.Cells.Item(20, 20).Value = 0.003974592724481738
$ws.Cells.Item(21, 7).Value = 2.891468666666666
$ws.Cells.Item(21, 8).Value = 8.674405999999999
$ws.Cells.Item(21, 9).Value = 0.1382948604962879
$ws.Cells.Item(21, 10).Value = 0.1382948604962879
$ws.Cells.Item(21, 13).Value = 61.75795633333333
$ws.Cells.Item(21, 14).Value = 185.273869
$ws.Cells.Item(21, 15).Value = 0.2272725152794058
$ws.Cells.Item(21, 16).Value = 0.2272725152794058
$ws.Cells.Item(21, 17).Value = 178.5711956552015
$ws.Cells.Item(21, 18).Value = 1607.140760896814
$ws.Cells.Item(21, 19).Value = 0.03143062079520589
$ws.Cells.Item(21, 20).Value = 0.03143062079520589
$ws.Cells.Item(22, 7).Value = 6.741071666666667
$ws.Cells.Item(22, 8).Value = 20.223215
$ws.Cells.Item(22, 9).Value = 0.3224159322507428
$ws.Cells.Item(22, 10).Value = 0.3224159322507428
$ws.Cells.Item(22, 13).Value = 86.89540866666668
$ws.Cells.Item(22, 14).Value = 260.686226
$ws.Cells.Item(22, 15).Value = 0.319779657009892
$ws.Cells.Item(22, 16).Value = 0.3197796570098919
$ws.Cells.Item(22, 17).Value = 585.7681773262879
$ws.Cells.Item(22, 18).Value = 5271.91359593659
$ws.Cells.Item(22, 19).Value = 0.1031020562296671
$ws.Cells.Item(22, 20).Value = 0.1031020562296671
$ws.Cells.Item(23, 7).Value = 6.741071666666667
$ws.Cells.Item(23, 8).Value = 20.223215
$ws.Cells.Item(23, 9).Value = 0.3224159322507428
$ws.Cells.Item(23, 10).Value = 0.3224159322507428
$ws.Cells.Item(23, 15).Value = 0.1999969065479545
$ws.Cells.Item(23, 16).Value = 0.1999969065479545
$ws.Cells.Item(23, 17).Value = 366.3517076568361
$ws.Cells.Item(23, 18).Value = 3297.165368911525
$ws.Cells.Item(23, 19).Value = 0.06448218907192346
$ws.Cells.Item(23, 20).Value = 0.06448218907192345
$ws.Cells.Item(24, 7).Value = 6.741071666666667
$ws.Cells.Item(24, 8).Value = 20.223215
$ws.Cells.Item(24, 9).Value = 0.3224159322507428
$ws.Cells.Item(24, 10).Value = 0.3224159322507428
$ws.Cells.Item(24, 13).Value = 60.92601633333334
$ws.Cells.Item(24, 14).Value = 182.778049
$ws.Cells.Item(24, 15).Value = 0.224210932487692
$ws.Cells.Item(24, 16).Value = 0.224210932487692
$ws.Cells.Item(24, 17).Value = 410.7066424675039
$ws.Cells.Item(24, 18).Value = 3696.359782207535
$ws.Cells.Item(24, 19).Value = 0.07228917681882757
$ws.Cells.Item(24, 20).Value = 0.07228917681882757
$ws.Cells.Item(25, 7).Value = 6.741071666666667
$ws.Cells.Item(25, 8).Value = 20.223215
$ws.Cells.Item(25, 9).Value = 0.3224159322507428
$ws.Cells.Item(25, 10).Value = 0.3224159322507428
$ws.Cells.Item(25, 13).Value = 7.809668333333332
$ws.Cells.Item(25, 14).Value = 23.429005
$ws.Cells.Item(25, 15).Value = 0.02873998867505581
$ws.Cells.Item(25, 16).Value = 0.02873998867505581
$ws.Cells.Item(25, 17).Value = 52.64553392789721
$ws.Cells.Item(25, 18).Value = 473.8098053510749
$ws.Cells.Item(25, 19).Value = 0.009266230241543912
$ws.Cells.Item(25, 20).Value = 0.009266230241543912
$ws.Cells.Item(26, 7).Value = 6.741071666666667
$ws.Cells.Item(26, 8).Value = 20.223215
$ws.Cells.Item(26, 9).Value = 0.3224159322507428
$ws.Cells.Item(26, 10).Value = 0.3224159322507428
$ws.Cells.Item(26, 13).Value = 61.75795633333333
$ws.Cells.Item(26, 14).Value = 185.273869
$ws.Cells.Item(26, 15).Value = 0.2272725152794058
$ws.Cells.Item(26, 16).Value = 0.2272725152794058
$ws.Cells.Item(26, 17).Value = 416.3148096298705
$ws.Cells.Item(26, 18).Value = 3746.833286668835
$ws.Cells.Item(26, 19).Value = 0.07327627988878083
$ws.Cells.Item(26, 20).Value = 0.07327627988878081
